$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.147.91'
$ws.Range("E2").Value = '  -1.97%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.85'
$ws.Range("E3").Value = '  -1.82%  '
$ws.Range("E4").Value = '  +0.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.09'
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5215'
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2665'
$ws.Range("E8").Value = '  -0.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06319'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.06'
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07746'
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.433'
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.650.71'
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.886.41'
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5474'
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8227'
$ws.Range("E16").Value = '  -3.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.94'
$ws.Range("E17").Value = '  -2.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.219.95'
$ws.Range("E18").Value = '  -1.88%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.666'
$ws.Range("E20").Value = '  -3.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.00'
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.16'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.099'
$ws.Range("E23").Value = '  -4.79%  '
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.33'
$ws.Range("E25").Value = '  -3.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1245'
$ws.Range("E26").Value = '  -2.69%  '
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.17'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.426'
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06014'
$ws.Range("E30").Value = '  -2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.283'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.572'
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.339'
$ws.Range("E33").Value = '  -3.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.650'
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9810'
$ws.Range("E35").Value = '  -3.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.412'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5916'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01594'
$ws.Range("E39").Value = '  -3.69%  '
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8640'
$ws.Range("E41").Value = '  -0.39%  '
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.038.87'
$ws.Range("E43").Value = '  -3.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.71'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.800.16'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.23'
$ws.Range("E46").Value = '  -0.47%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈108'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.120'
$ws.Range("E49").Value = '  -1.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05181'
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.467'
$ws.Range("E51").Value = '  +3.15%  '
